$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Diversion"
$ws.Range("A6").Value = "Carro"
$ws.Range("A7").Value = "Freelance"
